# Remove the leading "🧪 " (test-tube emoji + following space) that
# precedes the "Testing Checklist" heading. The remaining "Testing
# Checklist" run (and its formatting) is left completely untouched.
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("🧪 ", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 2)

Write-Host "Removed emoji prefix:" $found
